# Apply the localization-status.xlsx "Generate Report for Handoff" update.
#
# Summary of the change (derived from the OOXML diff):
#  1. In the source workbook, the Overview sheet's column G ("Latest HO
#     Xliff Generate Date") and the de-de sheet's column H ("Latest Handoff
#     Datetime") happen to point at the very same shared string. The diff
#     only edits that shared string's text (2016-08-20 18:34:37 ->
#     2016-08-20 18:35:06), so BOTH locations - every data row (2-14) of
#     Overview!G and every data row (2-14) of de-de!H - pick up the new
#     timestamp together.
#  2. zh-cn sheet, column H ("Latest Handoff Datetime") has its own,
#     separate shared string; every data row (2-14) changes from
#     2016-08-20 18:34:31 -> 2016-08-20 18:34:57.
#  3. zh-cn sheet, column E ("Priority") for every data row (2-14) changes
#     from blank -> "ht" (a brand-new shared string is introduced).
#  4. de-de sheet, column E ("Priority") for every data row (2-14) changes
#     from blank -> "ht" as well.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2:G14").Value = "2016-08-20 18:35:06"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2:H14").Value = "2016-08-20 18:34:57"
$zhcn.Range("E2:E14").Value = "ht"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2:H14").Value = "2016-08-20 18:35:06"
$dede.Range("E2:E14").Value = "ht"
